$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Gimmick sheet: add a new "hitPoint" (int) column (column D)
# ---------------------------------------------------------------------------
$gim = $wb.Worksheets.Item("Gimmick")

# Copy formatting from column C (last existing data column) into new column D
$gim.Range("C1:C4").Copy()
$gim.Range("D1:D4").PasteSpecial(-4122)
$gim.Range("C5").Copy()
$gim.Range("D5").PasteSpecial(-4122)

# Match the column width of column D to column C (13.83203125)
$gim.Columns.Item(4).ColumnWidth = $gim.Columns.Item(3).ColumnWidth

# Header / schema rows
$gim.Range("D1").Value = "hitPoint"
$gim.Range("D4").Value = "int"

# Sample data row
$gim.Range("D5").Value = 5

# Restore a sensible selection on this sheet
[void]$gim.Range("E8").Select()

# ---------------------------------------------------------------------------
# Character sheet: add a new "power" (int) column (column G)
# ---------------------------------------------------------------------------
$chr = $wb.Worksheets.Item("Character")

# Copy formatting from column F (last existing data column) into new column G
$chr.Range("F1:F4").Copy()
$chr.Range("G1:G4").PasteSpecial(-4122)
$chr.Range("F5").Copy()
$chr.Range("G5").PasteSpecial(-4122)

# Match the column width of column G to column F (14.5)
$chr.Columns.Item(7).ColumnWidth = $chr.Columns.Item(6).ColumnWidth

# Header / schema rows
$chr.Range("G1").Value = "power"
$chr.Range("G4").Value = "int"

# Sample data row
$chr.Range("G5").Value = 2

# Restore a sensible selection on this sheet
[void]$chr.Range("J17:J18").Select()
